$d = $word.ActiveDocument

# Step 1: remove _GoBack bookmark from paragraph 2
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Step 2: delete the empty paragraph between table1 and "Please add developer comments" (index 6)
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Delete()

# Step 3: resize table 1 + add bookmark
$t1 = $d.Tables.Item(1)
$t1.PreferredWidthType = 3
$t1.PreferredWidth = 603
$t1.Rows.LeftIndent = -67.75
$t1.Columns.Item(1).Width = 603
$cell1 = $t1.Cell(1,1)
$d.Bookmarks.Add("_GoBack", $cell1.Range)

# Step 4: resize table 2
$t2 = $d.Tables.Item(2)
$t2.PreferredWidthType = 3
$t2.PreferredWidth = 603
$t2.Rows.LeftIndent = -67.75
$t2.Columns.Item(1).Width = 603

Write-Host "Done"
